$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.737.38"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.908.11"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "312.41"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "0.5232"
$ws.Range("E7").Value = "  +7.96%  "
$ws.Range("D8").Value = "0.3784"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "0.07250"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "21.23"
$ws.Range("E10").Value = "  +3.66%  "
$ws.Range("D11").Value = "0.8980"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").Value = "0.07623"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "1.885.66"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "5.450"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "92.16"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "0.9993"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "0.000008724"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "0.9984"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").Value = "27.752.01"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "14.47"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "5.136"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "2.124.74"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "6.585"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "153.48"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("D27").Value = "2.163"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").Value = "18.30"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").Value = "114.70"
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").Value = "4.846"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").Value = "0.09017"
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").Value = "4.871"
$ws.Range("E32").Value = "  +5.23%  "
$ws.Range("D33").Value = "3.170"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("D35").Value = "0.7698"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("D36").Value = "2.625"
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("D37").Value = "0.02082"
$ws.Range("D38").Value = "3.068"
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("D39").Value = "1.090"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "0.5493"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").Value = "0.05279"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "6.647"
$ws.Range("E42").Value = "  -4.06%  "
$ws.Range("D43").Value = "113.78"
$ws.Range("E43").Value = "  +3.77%  "
$ws.Range("D44").Value = "8.487"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("D45").Value = "0.1508"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").Value = "0.4787"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "10.44"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").Value = "0.9993"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "1.616"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").Value = "66.48"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  -0.94%  "

$fmtRange.Style = "Normal"
